# SectorGroup.xlsx — reorder the codeforiati group/category name/code columns.
#
# The sheet has columns (1-indexed): A=code, B=name, C=status,
# D=codeforiati:group-name, E=codeforiati:category-name,
# F=codeforiati:group-code, G=codeforiati:category-code.
#
# The edit re-labels / re-orders these four columns so the header (and every
# data row) becomes: D=category-code, E=group-code, F=group-name,
# G=category-name. In other words, for every row the 4-tuple (D,E,F,G)
# goes from (a,b,c,d) to (d,c,a,b) — a single 4-cycle:
#   new D = old G
#   new E = old F
#   new F = old D
#   new G = old E
#
# We do this with real copy/paste (not by re-typing the literal values),
# so codes that look like numbers ("110", "111", ...) remain stored as text
# (shared strings) exactly like they were before, instead of being coerced
# into numeric cells with a new number-format style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$srcD = "D1:D$lastRow"
$srcE = "E1:E$lastRow"
$srcF = "F1:F$lastRow"
$srcG = "G1:G$lastRow"
$scratch = "J1:J$lastRow"

$xlPasteValues = -4163

# 1) stash old D in a scratch column
$ws.Range($srcD).Copy()
$ws.Range($scratch).PasteSpecial($xlPasteValues)

# 2) D <- old G
$ws.Range($srcG).Copy()
$ws.Range($srcD).PasteSpecial($xlPasteValues)

# 3) G <- old E
$ws.Range($srcE).Copy()
$ws.Range($srcG).PasteSpecial($xlPasteValues)

# 4) E <- old F
$ws.Range($srcF).Copy()
$ws.Range($srcE).PasteSpecial($xlPasteValues)

# 5) F <- old D (from scratch)
$ws.Range($scratch).Copy()
$ws.Range($srcF).PasteSpecial($xlPasteValues)

# 6) clean up the scratch column
$ws.Range($scratch).ClearContents()
$excel.CutCopyMode = 0
